# Fill in the missing "Priority" value ("H") for requirement RFR07 in the
# "As a/an / I want to ... / so that ... / Priority / Status" user-story
# table. The target cell is currently an empty paragraph whose paragraph
# mark already carries bold + red (FF0000) run formatting; we add a run
# with text "H" using that same formatting, matching the other rows in
# the same table which already show "H" (or "M") as their priority.

$d = $word.ActiveDocument

$targetId = "RFR07"
$cr = [char]13
$bel = [char]7

$targetCell = $null

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Rows.Count -lt 1) { continue }

    # Locate the "Priority" column for this table by inspecting its header row.
    $priorityCol = -1
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        try {
            $header = $t.Cell(1, $c).Range.Text.TrimEnd($cr, $bel)
        } catch {
            continue
        }
        if ($header -eq "Priority") {
            $priorityCol = $c
        }
    }
    if ($priorityCol -eq -1) { continue }

    # Find the row whose first column is the target requirement ID and
    # whose Priority cell is still empty.
    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        try {
            $idText = $t.Cell($r, 1).Range.Text.TrimEnd($cr, $bel)
        } catch {
            continue
        }
        if ($idText -ne $targetId) { continue }

        try {
            $priorityCell = $t.Cell($r, $priorityCol)
            $priorityText = $priorityCell.Range.Text.TrimEnd($cr, $bel)
        } catch {
            continue
        }

        if ($priorityText -eq "") {
            $targetCell = $priorityCell
        }
    }
}

if ($targetCell -ne $null) {
    $targetCell.Range.Text = "H"
    $targetCell.Range.Font.Bold = 1
    $targetCell.Range.Font.Color = 255
}
